$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits -------------------------------------------------

# K6: "Data Extraction and Classification" -> "Data Extraction and Reverification"
$ws.Range("K6").Value = "Data Extraction and Reverification"

# D11: "Write Introduction" -> "Write Introduction, abstract"
$ws.Range("D11").Value = "Write Introduction, abstract"

# D17: "Complete reading all the research papers " ->
#      "Complete reading all the research papers and data reverification"
$ws.Range("D17").Value = "Complete reading all the research papers and data reverification"

# --- New cells -----------------------------------------------------------------

# G7: new Progress (%) value of 100% - copy number format from G3 (percent style)
# so it matches the rest of the column, then set the value.
$ws.Range("G3").Copy()
$ws.Range("G7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G7").Value = 1

# F11: new Supporting People value "Shifin" - copy format from E11 (matching text style)
$ws.Range("E11").Copy()
$ws.Range("F11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F11").Value = "Shifin"

$excel.CutCopyMode = $false

# --- Date / numeric updates ------------------------------------------------

$ws.Range("M6").Value = 45826   # 2025-06-18

$ws.Range("L7").Value = 45827   # 2025-06-19
$ws.Range("M7").Value = 45838   # 2025-06-30

$ws.Range("L8").Value = 45827   # 2025-06-19
$ws.Range("M8").Value = 45838   # 2025-06-30

$ws.Range("M9").Value = 45879   # 2025-08-10

$ws.Range("B17").Value = 45826  # 2025-06-18

$ws.Range("A18").Value = 45827  # 2025-06-19
$ws.Range("B18").Value = 45838  # 2025-06-30

$ws.Range("A19").Value = 45838  # 2025-06-30
$ws.Range("B19").Value = 45839  # 2025-07-01
